# Updates cryptos list prices/volumes (GitHub Actions refresh), per commit diff.
# Values are written with a leading apostrophe to force literal-text storage,
# matching the workbook's original inlineStr (text) cell type and preventing
# Excel from auto-converting numeric-looking strings (e.g. "1.000", "0.9999").
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2 - Bitcoin
$ws.Range("D2").Value = '''27.092.63'
$ws.Range("E2").Value = '''  -2.95%  '

# Row 3 - Ethereum
$ws.Range("D3").Value = '''1.718.38'
$ws.Range("E3").Value = '''  -2.85%  '

# Row 4 - TetherUSD
$ws.Range("D4").Value = '''1.000'
$ws.Range("E4").Value = '''  -0.29%  '

# Row 5 - BNB
$ws.Range("D5").Value = '''311.60'
$ws.Range("E5").Value = '''  -5.13%  '

# Row 6 - USDC
$ws.Range("D6").Value = '''1.000'
$ws.Range("E6").Value = '''  -0.16%  '

# Row 7 - XRP
$ws.Range("D7").Value = '''0.4793'
$ws.Range("E7").Value = '''  +7.56%  '

# Row 8 - Cardano
$ws.Range("D8").Value = '''0.3471'
$ws.Range("E8").Value = '''  -1.26%  '

# Row 9 - OKB
$ws.Range("D9").Value = '''42.44'
$ws.Range("E9").Value = '''  +1.58%  '

# Row 10 - Dogecoin
$ws.Range("D10").Value = '''0.07272'
$ws.Range("E10").Value = '''  -1.46%  '

# Row 11 - Polygon
$ws.Range("D11").Value = '''1.043'
$ws.Range("E11").Value = '''  -4.94%  '

# Row 12 - BinanceUSD
$ws.Range("D12").Value = '''1.0000'
$ws.Range("E12").Value = '''  -0.29%  '

# Row 13 - Solana
$ws.Range("E13").Value = '''  -4.47%  '

# Row 14 - Polkadot
$ws.Range("D14").Value = '''5.854'
$ws.Range("E14").Value = '''  -2.73%  '

# Row 15 - WrappedEther
$ws.Range("D15").Value = '''1.699.86'
$ws.Range("E15").Value = '''  -4.05%  '

# Row 16 - Chainlink
$ws.Range("D16").Value = '''6.844'
$ws.Range("E16").Value = '''  -4.73%  '

# Row 17 - Litecoin
$ws.Range("D17").Value = '''88.00'
$ws.Range("E17").Value = '''  -4.95%  '

# Row 18 - ShibaInu
$ws.Range("D18").Value = '''0.00001035'
$ws.Range("E18").Value = '''  -2.44%  '

# Row 19 - TRON
$ws.Range("D19").Value = '''0.06406'
$ws.Range("E19").Value = '''  -0.18%  '

# Row 20 - Dai
$ws.Range("D20").Value = '''0.9999'
$ws.Range("E20").Value = '''  -0.20%  '

# Row 21 - Avalanche
$ws.Range("D21").Value = '''16.47'
$ws.Range("E21").Value = '''  -3.02%  '

# Row 22 - Uniswap
$ws.Range("D22").Value = '''5.631'
$ws.Range("E22").Value = '''  -2.30%  '

# Row 23 - WrappedBTC
$ws.Range("D23").Value = '''27.148.11'
$ws.Range("E23").Value = '''  -2.87%  '

# Row 24 - Cosmos
$ws.Range("D24").Value = '''10.81'
$ws.Range("E24").Value = '''  -3.61%  '

# Row 25 - Toncoin
$ws.Range("E25").Value = '''  -0.63%  '

# Row 26 - EthereumClassic (swapped into row 26)
$ws.Range("B26").Value = '''EthereumClassic'
$ws.Range("C26").Value = '''https://coinranking.com/coin/hnfQfsYfeIGUQ+ethereumclassic-etc'
$ws.Range("D26").Value = '''20.07'
$ws.Range("E26").Value = '''  -0.54%  '

# Row 27 - Monero (swapped into row 27)
$ws.Range("B27").Value = '''Monero'
$ws.Range("C27").Value = '''https://coinranking.com/coin/3mVx2FX_iJFp5+monero-xmr'
$ws.Range("D27").Value = '''150.47'
$ws.Range("E27").Value = '''  -6.54%  '

# Row 28 - WrappedliquidstakedEther2.0
$ws.Range("D28").Value = '''1.919.43'
$ws.Range("E28").Value = '''  -2.86%  '

# Row 29 - LidoDAOToken
$ws.Range("D29").Value = '''2.081'
$ws.Range("E29").Value = '''  -2.74%  '

# Row 30 - BitcoinCash
$ws.Range("D30").Value = '''121.13'

# Row 31 - ImmutableX
$ws.Range("D31").Value = '''1.034'
$ws.Range("E31").Value = '''  -4.80%  '

# Row 32 - Stellar
$ws.Range("D32").Value = '''0.09226'
$ws.Range("E32").Value = '''  +0.95%  '

# Row 33 - HuobiToken
$ws.Range("D33").Value = '''3.593'
$ws.Range("E33").Value = '''  -2.29%  '

# Row 34 - Filecoin
$ws.Range("D34").Value = '''5.342'
$ws.Range("E34").Value = '''  -5.04%  '

# Row 35 - WEMIXTOKEN
$ws.Range("D35").Value = '''1.469'
$ws.Range("E35").Value = '''  +6.45%  '

# Row 36 - VeChain
$ws.Range("D36").Value = '''0.02191'
$ws.Range("E36").Value = '''  -3.91%  '

# Row 37 - Hedera
$ws.Range("D37").Value = '''0.05885'
$ws.Range("E37").Value = '''  -3.67%  '

# Row 38 - Aptos
$ws.Range("D38").Value = '''10.99'
$ws.Range("E38").Value = '''  -7.08%  '

# Row 39 - Algorand
$ws.Range("D39").Value = '''0.1992'
$ws.Range("E39").Value = '''  -4.33%  '

# Row 40 - Frax
$ws.Range("D40").Value = '''0.9994'
$ws.Range("E40").Value = '''  -0.26%  '

# Row 41 - InternetComputer(DFINITY)
$ws.Range("D41").Value = '''4.727'
$ws.Range("E41").Value = '''  -4.61%  '

# Row 42 - TheSandbox
$ws.Range("D42").Value = '''0.5961'
$ws.Range("E42").Value = '''  -4.76%  '

# Row 43 - TrustWalletToken
$ws.Range("D43").Value = '''1.088'
$ws.Range("E43").Value = '''  -7.78%  '

# Row 44 - FraxShare
$ws.Range("D44").Value = '''7.485'
$ws.Range("E44").Value = '''  -4.22%  '

# Row 45 - EnergySwap
$ws.Range("D45").Value = '''12.84'
$ws.Range("E45").Value = '''  -3.20%  '

# Row 46 - PancakeSwap
$ws.Range("D46").Value = '''3.594'
$ws.Range("E46").Value = '''  -3.71%  '

# Row 47 - Decentraland
$ws.Range("D47").Value = '''0.5599'
$ws.Range("E47").Value = '''  -4.09%  '

# Row 48 - Quant
$ws.Range("D48").Value = '''118.61'
$ws.Range("E48").Value = '''  -3.28%  '

# Row 49 - NEARProtocol
$ws.Range("D49").Value = '''1.833'
$ws.Range("E49").Value = '''  -5.53%  '

# Row 50 - Cronos
$ws.Range("D50").Value = '''0.06703'
$ws.Range("E50").Value = '''  -2.38%  '

# Row 51 - EOS
$ws.Range("D51").Value = '''1.089'
$ws.Range("E51").Value = '''  -3.94%  '
